$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 204.558136
$ws.Cells.Item(2, 8).Value = 613.674408
$ws.Cells.Item(2, 9).Value = 0.60178627893129
$ws.Cells.Item(2, 10).Value = 0.6017862789312901
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1451143333333333
$ws.Cells.Item(2, 14).Value = 0.435343
$ws.Cells.Item(2, 15).Value = 0.140827770705994
$ws.Cells.Item(2, 16).Value = 0.1408277707059941
$ws.Cells.Item(2, 17).Value = 29.68431753354933
$ws.Cells.Item(2, 18).Value = 267.158857801944
$ws.Cells.Item(2, 19).Value = 0.08474822010334908
$ws.Cells.Item(2, 20).Value = 0.08474822010334911
$ws.Cells.Item(3, 7).Value = 204.558136
$ws.Cells.Item(3, 8).Value = 613.674408
$ws.Cells.Item(3, 9).Value = 0.60178627893129
$ws.Cells.Item(3, 10).Value = 0.6017862789312901
$ws.Cells.Item(3, 15).Value = 0.0422072807203407
$ws.Cells.Item(3, 16).Value = 0.0422072807203407
$ws.Cells.Item(3, 17).Value = 8.896642450912001
$ws.Cells.Item(3, 18).Value = 80.06978205820801
$ws.Cells.Item(3, 19).Value = 0.0253997624085022
$ws.Cells.Item(3, 20).Value = 0.02539976240850221
$ws.Cells.Item(4, 7).Value = 204.558136
$ws.Cells.Item(4, 8).Value = 613.674408
$ws.Cells.Item(4, 9).Value = 0.60178627893129
$ws.Cells.Item(4, 10).Value = 0.6017862789312901
$ws.Cells.Item(4, 13).Value = 0.841832
$ws.Cells.Item(4, 14).Value = 2.525496
$ws.Cells.Item(4, 15).Value = 0.8169649485736653
$ws.Cells.Item(4, 16).Value = 0.8169649485736653
$ws.Cells.Item(4, 17).Value = 172.203584745152
$ws.Cells.Item(4, 18).Value = 1549.832262706368
$ws.Cells.Item(4, 19).Value = 0.4916382964194387
$ws.Cells.Item(4, 20).Value = 0.4916382964194388
$ws.Cells.Item(5, 9).Value = 0.3090998990957371
$ws.Cells.Item(5, 10).Value = 0.3090998990957372
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.1451143333333333
$ws.Cells.Item(5, 14).Value = 0.435343
$ws.Cells.Item(5, 15).Value = 0.140827770705994
$ws.Cells.Item(5, 16).Value = 0.1408277707059941
$ws.Cells.Item(5, 17).Value = 15.24697367749978
$ws.Cells.Item(5, 18).Value = 137.222763097498
$ws.Cells.Item(5, 19).Value = 0.04352984971510036
$ws.Cells.Item(5, 20).Value = 0.04352984971510037
$ws.Cells.Item(6, 9).Value = 0.3090998990957371
$ws.Cells.Item(6, 10).Value = 0.3090998990957372
$ws.Cells.Item(6, 15).Value = 0.0422072807203407
$ws.Cells.Item(6, 16).Value = 0.0422072807203407
$ws.Cells.Item(6, 18).Value = 41.12682927693601
$ws.Cells.Item(6, 19).Value = 0.01304626621176276
$ws.Cells.Item(6, 20).Value = 0.01304626621176276
$ws.Cells.Item(7, 9).Value = 0.3090998990957371
$ws.Cells.Item(7, 10).Value = 0.3090998990957372
$ws.Cells.Item(7, 13).Value = 0.841832
$ws.Cells.Item(7, 14).Value = 2.525496
$ws.Cells.Item(7, 15).Value = 0.8169649485736653
$ws.Cells.Item(7, 16).Value = 0.8169649485736653
$ws.Cells.Item(7, 17).Value = 88.45018992985068
$ws.Cells.Item(7, 18).Value = 796.051709368656
$ws.Cells.Item(7, 19).Value = 0.252523783168874
$ws.Cells.Item(7, 20).Value = 0.252523783168874
$ws.Cells.Item(8, 7).Value = 0.1651866666666667
$ws.Cells.Item(8, 8).Value = 0.49556
$ws.Cells.Item(8, 9).Value = 0.0004859599887163456
$ws.Cells.Item(8, 10).Value = 0.0004859599887163457
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1451143333333333
$ws.Cells.Item(8, 14).Value = 0.435343
$ws.Cells.Item(8, 15).Value = 0.140827770705994
$ws.Cells.Item(8, 16).Value = 0.1408277707059941
$ws.Cells.Item(8, 17).Value = 0.02397095300888889
$ws.Cells.Item(8, 18).Value = 0.21573857708
$ws.Cells.Item(8, 19).Value = 0.00006843666186323297
$ws.Cells.Item(8, 20).Value = 0.00006843666186323298
$ws.Cells.Item(9, 7).Value = 0.1651866666666667
$ws.Cells.Item(9, 8).Value = 0.49556
$ws.Cells.Item(9, 9).Value = 0.0004859599887163456
$ws.Cells.Item(9, 10).Value = 0.0004859599887163457
$ws.Cells.Item(9, 15).Value = 0.0422072807203407
$ws.Cells.Item(9, 16).Value = 0.0422072807203407
$ws.Cells.Item(9, 17).Value = 0.007184298506666667
$ws.Cells.Item(9, 18).Value = 0.06465868656
$ws.Cells.Item(9, 19).Value = 0.0000205110496626044
$ws.Cells.Item(9, 20).Value = 0.0000205110496626044
$ws.Cells.Item(10, 7).Value = 0.1651866666666667
$ws.Cells.Item(10, 8).Value = 0.49556
$ws.Cells.Item(10, 9).Value = 0.0004859599887163456
$ws.Cells.Item(10, 10).Value = 0.0004859599887163457
$ws.Cells.Item(10, 13).Value = 0.841832
$ws.Cells.Item(10, 14).Value = 2.525496
$ws.Cells.Item(10, 15).Value = 0.8169649485736653
$ws.Cells.Item(10, 16).Value = 0.8169649485736653
$ws.Cells.Item(10, 17).Value = 0.1390594219733333
$ws.Cells.Item(10, 18).Value = 1.25153479776
$ws.Cells.Item(10, 19).Value = 0.0003970122771905083
$ws.Cells.Item(10, 20).Value = 0.0003970122771905083
$ws.Cells.Item(11, 7).Value = 29.46746633333333
$ws.Cells.Item(11, 8).Value = 88.402399
$ws.Cells.Item(11, 9).Value = 0.08668986363011115
$ws.Cells.Item(11, 10).Value = 0.08668986363011116
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.1451143333333333
$ws.Cells.Item(11, 14).Value = 0.435343
$ws.Cells.Item(11, 15).Value = 0.140827770705994
$ws.Cells.Item(11, 16).Value = 0.1408277707059941
$ws.Cells.Item(11, 17).Value = 4.276151731984111
$ws.Cells.Item(11, 18).Value = 38.485365587857
$ws.Cells.Item(11, 19).Value = 0.01220834023783518
$ws.Cells.Item(11, 20).Value = 0.01220834023783519
$ws.Cells.Item(12, 7).Value = 29.46746633333333
$ws.Cells.Item(12, 8).Value = 88.402399
$ws.Cells.Item(12, 9).Value = 0.08668986363011115
$ws.Cells.Item(12, 10).Value = 0.08668986363011116
$ws.Cells.Item(12, 15).Value = 0.0422072807203407
$ws.Cells.Item(12, 16).Value = 0.0422072807203407
$ws.Cells.Item(12, 17).Value = 1.281599045769333
$ws.Cells.Item(12, 18).Value = 11.534391411924
$ws.Cells.Item(12, 19).Value = 0.003658943409844154
$ws.Cells.Item(12, 20).Value = 0.003658943409844156
$ws.Cells.Item(13, 7).Value = 29.46746633333333
$ws.Cells.Item(13, 8).Value = 88.402399
$ws.Cells.Item(13, 9).Value = 0.08668986363011115
$ws.Cells.Item(13, 10).Value = 0.08668986363011116
$ws.Cells.Item(13, 13).Value = 0.841832
$ws.Cells.Item(13, 14).Value = 2.525496
$ws.Cells.Item(13, 15).Value = 0.8169649485736653
$ws.Cells.Item(13, 16).Value = 0.8169649485736653
$ws.Cells.Item(13, 17).Value = 24.80665611832267
$ws.Cells.Item(13, 18).Value = 223.259905064904
$ws.Cells.Item(13, 19).Value = 0.07082257998243181
$ws.Cells.Item(13, 20).Value = 0.07082257998243183
$ws.Cells.Item(14, 7).Value = 0.54608
$ws.Cells.Item(14, 8).Value = 1.63824
$ws.Cells.Item(14, 9).Value = 0.001606503938805929
$ws.Cells.Item(14, 10).Value = 0.001606503938805929
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.1451143333333333
$ws.Cells.Item(14, 14).Value = 0.435343
$ws.Cells.Item(14, 15).Value = 0.140827770705994
$ws.Cells.Item(14, 16).Value = 0.1408277707059941
$ws.Cells.Item(14, 17).Value = 0.07924403514666666
$ws.Cells.Item(14, 18).Value = 0.7131963163199999
$ws.Cells.Item(14, 19).Value = 0.0002262403683324376
$ws.Cells.Item(14, 20).Value = 0.0002262403683324377
$ws.Cells.Item(15, 7).Value = 0.54608
$ws.Cells.Item(15, 8).Value = 1.63824
$ws.Cells.Item(15, 9).Value = 0.001606503938805929
$ws.Cells.Item(15, 10).Value = 0.001606503938805929
$ws.Cells.Item(15, 15).Value = 0.0422072807203407
$ws.Cells.Item(15, 16).Value = 0.0422072807203407
$ws.Cells.Item(15, 17).Value = 0.02375011136
$ws.Cells.Item(15, 18).Value = 0.21375100224
$ws.Cells.Item(15, 19).Value = 0.00006780616272351486
$ws.Cells.Item(15, 20).Value = 0.00006780616272351489
$ws.Cells.Item(16, 7).Value = 0.54608
$ws.Cells.Item(16, 8).Value = 1.63824
$ws.Cells.Item(16, 9).Value = 0.001606503938805929
$ws.Cells.Item(16, 10).Value = 0.001606503938805929
$ws.Cells.Item(16, 13).Value = 0.841832
$ws.Cells.Item(16, 14).Value = 2.525496
$ws.Cells.Item(16, 15).Value = 0.8169649485736653
$ws.Cells.Item(16, 16).Value = 0.8169649485736653
$ws.Cells.Item(16, 17).Value = 0.45970761856
$ws.Cells.Item(16, 18).Value = 4.137368567039999
$ws.Cells.Item(16, 19).Value = 0.001312457407749976
$ws.Cells.Item(16, 20).Value = 0.001312457407749976
$ws.Cells.Item(17, 7).Value = 0.112681
$ws.Cells.Item(17, 8).Value = 0.338043
$ws.Cells.Item(17, 9).Value = 0.000331494415339494
$ws.Cells.Item(17, 10).Value = 0.000331494415339494
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.1451143333333333
$ws.Cells.Item(17, 14).Value = 0.435343
$ws.Cells.Item(17, 15).Value = 0.140827770705994
$ws.Cells.Item(17, 16).Value = 0.1408277707059941
$ws.Cells.Item(17, 17).Value = 0.01635162819433333
$ws.Cells.Item(17, 18).Value = 0.147164653749
$ws.Cells.Item(17, 19).Value = 0.00004668361951374781
$ws.Cells.Item(17, 20).Value = 0.00004668361951374782
$ws.Cells.Item(18, 7).Value = 0.112681
$ws.Cells.Item(18, 8).Value = 0.338043
$ws.Cells.Item(18, 9).Value = 0.000331494415339494
$ws.Cells.Item(18, 10).Value = 0.000331494415339494
$ws.Cells.Item(18, 15).Value = 0.0422072807203407
$ws.Cells.Item(18, 16).Value = 0.0422072807203407
$ws.Cells.Item(18, 17).Value = 0.004900722052000001
$ws.Cells.Item(18, 18).Value = 0.04410649846800001
$ws.Cells.Item(18, 19).Value = 0.00001399147784545923
$ws.Cells.Item(18, 20).Value = 0.00001399147784545924
$ws.Cells.Item(19, 7).Value = 0.112681
$ws.Cells.Item(19, 8).Value = 0.338043
$ws.Cells.Item(19, 9).Value = 0.000331494415339494
$ws.Cells.Item(19, 10).Value = 0.000331494415339494
$ws.Cells.Item(19, 13).Value = 0.841832
$ws.Cells.Item(19, 14).Value = 2.525496
$ws.Cells.Item(19, 15).Value = 0.8169649485736653
$ws.Cells.Item(19, 16).Value = 0.8169649485736653
$ws.Cells.Item(19, 17).Value = 0.09485847159200002
$ws.Cells.Item(19, 18).Value = 0.853726244328
$ws.Cells.Item(19, 19).Value = 0.000270819317980287
$ws.Cells.Item(19, 20).Value = 0.000270819317980287
